# Results for Modulo IK
# The "results (2)" sheet's query-table data (A1:F37) was re-sorted in
# Excel by Column1 (B), then Column2 (C), then Column3 (D), ascending,
# and the active selection ended up on B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results (2)")
$ws.Activate()

$rng = $ws.Range("A1:F37")
$rng.Sort($ws.Range("B1"), 1, $ws.Range("C1"), $null, 1, $ws.Range("D1"), 1, 1, $false, $null, $null, 1)

$ws.Range("B8").Select()
